$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for the new columns ---------------------------------------
# Old layout:  A=ID  B=Email  C=EmailPassword  D=MLBPassword  E=7-8(results)
# New layout:  A=ID  B=Email  C=EmailPassword  D=MLBPassword
#              E=Strategy  F=VM  G=7-8(results, shifted)  H=7-9(results)
# Inserting whole columns (rather than just writing into blank cells) makes
# Excel carry the existing header formatting (bold/border/center, style
# index 1) into the new header cells automatically, and shifts the old
# "7-8" results column from E to G intact.
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(8).Insert()

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("E1").Value = "Strategy"
$ws.Range("F1").Value = "VM"
$ws.Range("H1").Value = "7-9"

# --- Row 2 ------------------------------------------------------------------
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 1
$ws.Range("H2").Value = "Done. 1: ('Robinson', 'Cano', 'sea'), 2: ()"

# --- Row 3 ------------------------------------------------------------------
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 2
$ws.Range("H3").Value = "NOT DONE"
